$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "AddCustomerTest"

$ws.Range("A1").Value = "firstname"
$ws.Range("B1").Value = "lastname"
$ws.Range("A2").Value = "Dan"
$ws.Range("B2").Value = "Car"

$ws.Range("C1").Value = "postcode"
$ws.Range("C2").Value = "syd2148"

$ws.Range("D1").Value = "alerttest"
$ws.Range("D2").Value = "Customer added successfully"

$ws.Range("D2").Select()
